$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source row " basal diameter" / QUALITY / 1 (row 433) was removed entirely,
# and all subsequent rows (434:722) shift up by one row to fill the gap.
$ws.Rows.Item(433).Delete()
